# Ozurgeti disability_prevalence.xlsx update
# Retitle the header, add a new "family with disabilities" data row above
# the existing (renamed) data row, and tidy up row heights / column width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert a new row above the current data row (row 4). This pushes the
#    old "Number of disability persons" row to 5 and the Source row to 6,
#    carrying their values/formats/merge down with them.
# ---------------------------------------------------------------------
$ws.Rows.Item(4).Insert(-4121)

# ---------------------------------------------------------------------
# 2. Title row (row 1): new wording, merged across A1:I1, bold Arial 11,
#    centered + wrapped, taller row.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Ozurgeti Municipality"
$ws.Range("A1:I1").Font.Name = "Arial"
$ws.Range("A1:I1").Font.Size = 11
$ws.Range("A1:I1").Font.Bold = $true
$ws.Range("A1:I1").HorizontalAlignment = -4108
$ws.Range("A1:I1").VerticalAlignment = -4108
$ws.Range("A1:I1").WrapText = $true
$ws.Range("A1:I1").Merge()
$ws.Rows.Item(1).RowHeight = 51

# ---------------------------------------------------------------------
# 3. Row 3 / cell A3: swap in the "Sylfaen" font + top border (the year
#    header cells B3:I3 keep their existing formatting untouched).
# ---------------------------------------------------------------------
$ws.Range("A3").Font.Name = "Sylfaen"
$ws.Range("A3").Font.Size = 11
$ws.Range("A3").Font.Color = 0
$ws.Range("A3").Borders.Item(9).LineStyle = 1

# ---------------------------------------------------------------------
# 4. Row 4 (new): "family with disabilities Persons " + historical values.
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "family with disabilities Persons "
$ws.Range("A4").Font.Name = "Arial"
$ws.Range("A4").Font.Size = 10
$ws.Range("A4").Font.Color = 0
$ws.Range("A4").Interior.Color = 16777215
$ws.Range("A4").Borders.Item(9).LineStyle = 1
$ws.Range("A4").HorizontalAlignment = -4131
$ws.Range("A4").VerticalAlignment = -4108
$ws.Range("A4").WrapText = $true

$ws.Cells.Item(4,2).Value = 1453
$ws.Cells.Item(4,3).Value = 1450
$ws.Cells.Item(4,4).Value = 1423
$ws.Cells.Item(4,5).Value = 1492
$ws.Cells.Item(4,6).Value = 1498
$ws.Cells.Item(4,7).Value = 1158
$ws.Cells.Item(4,8).Value = 1440
$ws.Cells.Item(4,9).Value = 1452

$ws.Range("B4:I4").Font.Name = "Arial"
$ws.Range("B4:I4").Font.Size = 10
$ws.Range("B4:I4").Font.Color = 0
$ws.Range("B4:I4").Interior.Color = 16777215
$ws.Range("B4:I4").NumberFormat = "#\ ##0"

$ws.Rows.Item(4).RowHeight = 24.75

# ---------------------------------------------------------------------
# 5. Row 5 (was row 4): rename + replace the values, keep the same
#    left/center/wrap styling but move the thin border to the bottom.
# ---------------------------------------------------------------------
$ws.Range("A5").Value = "disabilities Persons "
$ws.Range("A5").Borders.Item(8).LineStyle = 0
$ws.Range("A5").Borders.Item(9).LineStyle = 1

$ws.Cells.Item(5,2).Value = 1613
$ws.Cells.Item(5,3).Value = 1607
$ws.Cells.Item(5,4).Value = 1584
$ws.Cells.Item(5,5).Value = 1648
$ws.Cells.Item(5,6).Value = 1660
$ws.Cells.Item(5,7).Value = 1397
$ws.Cells.Item(5,8).Value = 1597
$ws.Cells.Item(5,9).Value = 1607

$ws.Range("B5:H5").Borders.Item(8).LineStyle = 0
$ws.Range("I5").Borders.Item(8).LineStyle = 0
$ws.Range("I5").Borders.Item(9).LineStyle = 1

$ws.Rows.Item(5).RowHeight = 21

# ---------------------------------------------------------------------
# 6. Row 6 (source row): drop the top border from A6, keep it on B6:H6.
# ---------------------------------------------------------------------
$ws.Range("A6").Borders.Item(8).LineStyle = 0
$ws.Range("B6:H6").Borders.Item(8).LineStyle = 1
$ws.Rows.Item(6).RowHeight = 27.75

# ---------------------------------------------------------------------
# 7. Column A width + selection cosmetics matching the authored file.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 19.92

$ws.Range("A1:I1").Select()
